# Apply the "cryptos list" update (Price + Volume(1h) columns) from the commit diff.
# (commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores every value as TEXT, even the plain-looking
# single-decimal ones (e.g. "591.20", "7.63"). Assigning such a string straight
# to .Value makes Excel auto-detect it as a number (and round-trip it as a
# float), so for just those cells we pin the format to Text first, then restore
# the default "Normal" style afterwards so no extra formatting is left behind.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.394.49"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").Value = "3.509.59"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "591.20"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("D6").Value = "134.58"
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "7.63"
$ws.Range("E9").Value = "  +6.09%  "
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("E11").Value = "  +2.85%  "
$ws.Range("D12").Value = "4.105.87"
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("E13").Value = "  +0.98%  "
$ws.Range("E14").Value = "  -0.64%  "
$ws.Range("D15").Value = "3.507.85"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").Value = "64.384.99"
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "25.70"
$ws.Range("E17").Value = "  +2.72%  "
$ws.Range("E18").Value = "  +0.17%  "
$ws.Range("E19").Value = "  +2.02%  "
$ws.Range("D20").Value = "13.63"
$ws.Range("E20").Value = "  -1.35%  "
$ws.Range("D21").Value = "394.75"
$ws.Range("E21").Value = "  +2.41%  "
$ws.Range("E22").Value = "  +1.40%  "
$ws.Range("D23").Value = "3.650.33"
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("D24").Value = "74.72"
$ws.Range("E24").Value = "  +0.76%  "
$ws.Range("E25").Value = "  +0.08%  "
$ws.Range("E26").Value = "  -0.07%  "
$ws.Range("D27").Value = "0.0000118"
$ws.Range("E27").Value = "  +3.05%  "
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("D29").Value = "7.40"
$ws.Range("E29").Value = "  -1.44%  "
$ws.Range("E30").Value = "  +1.05%  "
$ws.Range("D31").Value = "8.25"
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("E32").Value = "  -6.25%  "
$ws.Range("E33").Value = "  +5.85%  "
$ws.Range("D34").Value = "3.540.71"
$ws.Range("E34").Value = "  +0.46%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").Value = "23.36"
$ws.Range("E36").Value = "  -1.00%  "
$ws.Range("D37").Value = "5.36"
$ws.Range("E37").Value = "  +0.65%  "
$ws.Range("D38").Value = "6.97"
$ws.Range("E38").Value = "  +1.57%  "
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("D40").Value = "167.20"
$ws.Range("E40").Value = "  +2.21%  "
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("E42").Value = "  +0.32%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "4.44"
$ws.Range("E44").Value = "  +0.27%  "
$ws.Range("D45").Value = "24.90"
$ws.Range("E45").Value = "  -4.18%  "
$ws.Range("E46").Value = "  +0.55%  "
$ws.Range("E47").Value = "  -3.64%  "
$ws.Range("D48").Value = "6.81"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").Value = "2.380.41"
$ws.Range("E49").Value = "  -3.92%  "
$ws.Range("D50").Value = "0.899"
$ws.Range("E50").Value = "  -1.80%  "
$ws.Range("D51").Value = "0.0260"
$ws.Range("E51").Value = "  -0.27%  "

# Drop the temporary Text formatting back to Normal now that the text values
# are safely stored, so the cells keep the workbook default (no explicit style).
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
